$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Tasks and Bugs": update Progress % for the three DB-test related rows
# (rows 12-14, column E) from 0 to their new completion values.
$ws.Range("E12").Value = 100
$ws.Range("E13").Value = 10
$ws.Range("E14").Value = 10

# Move the viewport/selection to where the work is now happening:
# topLeftCell A7->A11, selection E21->E17.
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E17").Select()
